# Update "want to go" counts (column F) across the three sheets that carry
# this data: 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types,
# which mirrors the rows from the first two sheets).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 17
    7  = 1304
    8  = 511
    10 = 1287
    19 = 1656
    22 = 206
    23 = 2031
    26 = 921
    30 = 2811
    31 = 1604
    32 = 82
    34 = 651
    35 = 858
    36 = 1796
    38 = 1808
    39 = 201
    41 = 836
    43 = 852
    44 = 789
    45 = 1003
    46 = 62
    47 = 436
    48 = 3330
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 794

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 17
    8  = 1304
    9  = 511
    11 = 1287
    21 = 1656
    24 = 206
    25 = 2031
    29 = 2811
    30 = 1604
    31 = 82
    33 = 794
    35 = 651
    36 = 858
    37 = 1796
    40 = 1808
    41 = 836
    42 = 852
    43 = 789
    44 = 1003
    45 = 436
    48 = 3330
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}

$wb.Save()
